# Give the user a way to clear items from the cart: drop buyer orders #2-#4
# (rows 3-5), leaving only the header and the first buyer's order, and
# update that remaining order to a single-item cart.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra order rows (3,4,5) entirely, shifting rows up.
$ws.Range("A3:C5").EntireRow.Delete()

# The remaining cart (row 2) now only has the one item, at the lower total.
$ws.Range("B2").Value = "Apple Juice: 5pcs."
$ws.Range("C2").Value = "PHP 100.00"

# Reset the header/data cell formatting (bold fonts, colored fills) back to
# the plain default style now that the colour-coded cart rows are gone.
$ws.Range("A1:C2").ClearFormats()
